$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H18").Value = 1332.25
$ws.Range("J18").Value = 2998.5
$ws.Range("L18").Value = 2998.5
$ws.Range("N18").Value = -3566.5
$ws.Range("H40").Value = 3633.0833
$ws.Range("I40").Value = 4350
$ws.Range("J40").Value = 2199.25
$ws.Range("K40").Value = 4350
$ws.Range("L40").Value = 2199.25
$ws.Range("M40").Value = -4175
$ws.Range("N40").Value = -2549.25
$ws.Range("H52").Value = 476
$ws.Range("I52").Value = 476
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 1428
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -1268
$ws.Range("N52").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H86").Value = 1615.6666
$ws.Range("I86").Value = 1498
$ws.Range("K86").Value = 1498
$ws.Range("M86").Value = -375
$ws.Range("H89").Value = 1615.6666
$ws.Range("I89").Value = 1498
$ws.Range("K89").Value = 7490
$ws.Range("M89").Value = -1874
$ws.Range("H106").Value = 500005000
$ws.Range("I106").Value = 500005000
$ws.Range("K106").Value = 500005000
$ws.Range("M106").Value = -500004369

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 536.125
$ws.Range("I2").Value = 536.125
$ws.Range("K2").Value = 536.125
$ws.Range("M2").Value = -423.125
$ws.Range("H32").Value = 4764.294
$ws.Range("I32").Value = 4764.294
$ws.Range("K32").Value = 4764.294
$ws.Range("M32").Value = -4477.294
$ws.Range("H116").Value = 536.125
$ws.Range("I116").Value = 536.125
$ws.Range("K116").Value = 536.125
$ws.Range("M116").Value = 1757.875
$ws.Range("H122").Value = 700
$ws.Range("I122").Value = 700
$ws.Range("K122").Value = 2100
$ws.Range("M122").Value = 350
$ws.Range("H132").Value = 1341.5834
$ws.Range("I132").Value = 1440.2
$ws.Range("J132").Value = 848.5
$ws.Range("K132").Value = 4320.6
$ws.Range("L132").Value = 2545.5
$ws.Range("M132").Value = -1790.6
$ws.Range("N132").Value = -7605.5

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 536.125
$ws.Range("I3").Value = 536.125
$ws.Range("K3").Value = 536.125
$ws.Range("M3").Value = -422.125

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H7").Value = 89.42856999999999
$ws.Range("I7").Value = 132.88889
$ws.Range("J7").Value = 11.2
$ws.Range("K7").Value = 132.88889
$ws.Range("L7").Value = 11.2
$ws.Range("M7").Value = -19.88889
$ws.Range("N7").Value = -237.2
$ws.Range("H22").Value = 640.6786
$ws.Range("I22").Value = 845.38464
$ws.Range("J22").Value = 463.26666
$ws.Range("K22").Value = 845.38464
$ws.Range("L22").Value = 463.26666
$ws.Range("M22").Value = -495.38464
$ws.Range("N22").Value = -1163.26666
$ws.Range("H132").Value = 2245.111
$ws.Range("I132").Value = 2140.1428
$ws.Range("K132").Value = 6420.428400000001
$ws.Range("M132").Value = -3890.428400000001

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H3").Value = 9997.25
$ws.Range("I3").Value = 9996.333000000001
$ws.Range("K3").Value = 29988.999
$ws.Range("M3").Value = -29876.999
$ws.Range("H4").Value = 102881.836
$ws.Range("I4").Value = 43148.24
$ws.Range("J4").Value = 401549.8
$ws.Range("K4").Value = 129444.72
$ws.Range("L4").Value = 1204649.4
$ws.Range("M4").Value = -129332.72
$ws.Range("N4").Value = -1204873.4
$ws.Range("H33").Value = 405
$ws.Range("I33").Value = 395
$ws.Range("J33").Value = 431.66666
$ws.Range("K33").Value = 2370
$ws.Range("L33").Value = 2589.99996
$ws.Range("M33").Value = -2087
$ws.Range("N33").Value = -3155.99996
$ws.Range("H107").Value = 510.15384
$ws.Range("I107").Value = 237.6
$ws.Range("K107").Value = 712.8
$ws.Range("M107").Value = 1207.2
$ws.Range("H109").Value = 1488.5
$ws.Range("I109").Value = 1488.5
$ws.Range("K109").Value = 4465.5
$ws.Range("M109").Value = -3425.5
$ws.Range("H114").Value = 1550.875
$ws.Range("I114").Value = 1567.8334
$ws.Range("K114").Value = 4703.5002
$ws.Range("M114").Value = -1449.5002
$ws.Range("H129").Value = 1030
$ws.Range("I129").Value = 1030
$ws.Range("K129").Value = 3090
$ws.Range("M129").Value = 1910
$ws.Range("H131").Value = 4296.2
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 173
$ws.Range("I2").Value = 190.4
$ws.Range("K2").Value = 190.4
$ws.Range("M2").Value = -77.40000000000001
$ws.Range("H58").Value = 31020.5
$ws.Range("I58").Value = 31020.5
$ws.Range("K58").Value = 31020.5
$ws.Range("M58").Value = -30743.5
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H132").Value = 2679.3333
$ws.Range("I132").Value = 2012
$ws.Range("J132").Value = 4014
$ws.Range("K132").Value = 6036
$ws.Range("L132").Value = 12042
$ws.Range("M132").Value = -3506
$ws.Range("N132").Value = -17102
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1296.4445
$ws.Range("I16").Value = 1023.6
$ws.Range("K16").Value = 1023.6
$ws.Range("M16").Value = -853.6
$ws.Range("H17").Value = 500
$ws.Range("J17").Value = 500
$ws.Range("L17").Value = 500
$ws.Range("M17").Value = -840
$ws.Range("H22").Value = 932.6667
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 932.6667
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 932.6667
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1522.6667
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H27").Value = 932.6667
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 932.6667
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 932.6667
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1146.6667
$ws.Range("H40").Value = 1101.2
$ws.Range("I40").Value = 1101.2
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1101.2
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -965.2
$ws.Range("N40").ClearContents()
$ws.Range("H43").Value = 10006
$ws.Range("H55").Value = 740.6667
$ws.Range("I55").Value = 562.2
$ws.Range("K55").Value = 562.2
$ws.Range("M55").Value = -389.2
$ws.Range("H93").Value = 47622308
$ws.Range("I93").Value = 66669920
$ws.Range("K93").Value = 66669920
$ws.Range("M93").Value = -66668672
$ws.Range("H122").Value = 3771.8572
$ws.Range("I122").Value = 3450.75
$ws.Range("K122").Value = 10352.25
$ws.Range("M122").Value = -7902.25
$ws.Range("H132").Value = 3005.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3005.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 9016.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -14076.5

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 1199.8
$ws.Range("I122").Value = 999.6667
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 2999.0001
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -549.0001000000002
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
